$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I86").Value = 381119040
$ws.Range("J86").Value = 211112980
$ws.Range("K86").Value = 381119040
$ws.Range("L86").Value = 211112980
$ws.Range("M86").Value = -381117917
$ws.Range("N86").Value = -211115226
$ws.Range("I89").Value = 381119040
$ws.Range("J89").Value = 211112980
$ws.Range("K89").Value = 1905595200
$ws.Range("L89").Value = 1055564900
$ws.Range("M89").Value = -1905589584
$ws.Range("N89").Value = -1055576132
$ws.Range("H103").Value = 1169.3334
$ws.Range("I103").Value = 844.6667
$ws.Range("J103").Value = 1277.5555
$ws.Range("K103").Value = 2534.0001
$ws.Range("L103").Value = 3832.6665
$ws.Range("M103").Value = -1948.0001
$ws.Range("N103").Value = -5004.666499999999
$ws.Range("H106").Value = 3709.8462
$ws.Range("I106").Value = 3553.4
$ws.Range("K106").Value = 3553.4
$ws.Range("M106").Value = -2922.4
$ws.Range("H111").Value = 8933051
$ws.Range("I111").Value = 20835458
$ws.Range("J111").Value = 6246.875
$ws.Range("K111").Value = 62506374
$ws.Range("L111").Value = 18740.625
$ws.Range("M111").Value = -62503307
$ws.Range("N111").Value = -24874.625
$ws.Range("H113").Value = 27793424
$ws.Range("I113").Value = 2653.75
$ws.Range("J113").Value = 35733644
$ws.Range("K113").Value = 2653.75
$ws.Range("L113").Value = 35733644
$ws.Range("M113").Value = 600.25
$ws.Range("N113").Value = -35740152
$ws.Range("H118").Value = 815.5
$ws.Range("I118").Value = 815.5
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2446.5
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -789.5
$ws.Range("H126").Value = 58569.332
$ws.Range("I126").Value = 709
$ws.Range("J126").Value = 87499.5
$ws.Range("K126").Value = 709
$ws.Range("L126").Value = 87499.5
$ws.Range("M126").Value = 4231
$ws.Range("N126").Value = -97379.5
$ws.Range("H138").Value = 2525.09
$ws.Range("J138").Value = 2802.2026
$ws.Range("L138").Value = 8406.6078
$ws.Range("N138").Value = -18686.6078
$ws.Range("N118").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2670.818
$ws.Range("I14").Value = 3625
$ws.Range("J14").Value = 2125.5715
$ws.Range("K14").Value = 3625
$ws.Range("L14").Value = 2125.5715
$ws.Range("M14").Value = -3450
$ws.Range("N14").Value = -2475.5715
$ws.Range("H32").Value = 2156143
$ws.Range("I32").Value = 2179036
$ws.Range("K32").Value = 2179036
$ws.Range("M32").Value = -2178749
$ws.Range("H61").Value = 30308684
$ws.Range("I61").Value = 2556.65
$ws.Range("K61").Value = 2556.65
$ws.Range("M61").Value = -2344.65
$ws.Range("H122").Value = 3458.5186
$ws.Range("I122").Value = 2584.3125
$ws.Range("J122").Value = 4730.091
$ws.Range("K122").Value = 7752.9375
$ws.Range("L122").Value = 14190.273
$ws.Range("M122").Value = -5302.9375
$ws.Range("N122").Value = -19090.273
$ws.Range("H136").Value = 30308684
$ws.Range("I136").Value = 2556.65
$ws.Range("K136").Value = 7669.950000000001
$ws.Range("M136").Value = -5119.950000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5558063
$ws.Range("I99").Value = 1945.2142
$ws.Range("K99").Value = 1945.2142
$ws.Range("M99").Value = -447.2141999999999
$ws.Range("H105").Value = 2796.818
$ws.Range("I105").Value = 1395.7858
$ws.Range("K105").Value = 1395.7858
$ws.Range("M105").Value = 351.2141999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 624.8570999999999
$ws.Range("I22").Value = 683.1667
$ws.Range("K22").Value = 683.1667
$ws.Range("M22").Value = -333.1667
$ws.Range("H31").Value = 9585.767
$ws.Range("J31").Value = 12327.143
$ws.Range("L31").Value = 12327.143
$ws.Range("N31").Value = -12917.143
$ws.Range("H34").Value = 9585.767
$ws.Range("J34").Value = 12327.143
$ws.Range("L34").Value = 12327.143
$ws.Range("N34").Value = -12731.143
$ws.Range("H58").Value = 6490.8184
$ws.Range("I58").Value = 2995.6538
$ws.Range("K58").Value = 2995.6538
$ws.Range("M58").Value = -2792.6538
$ws.Range("H122").Value = 4390.8945
$ws.Range("I122").Value = 3488.8
$ws.Range("K122").Value = 10466.4
$ws.Range("M122").Value = -8016.400000000001
$ws.Range("H136").Value = 6490.8184
$ws.Range("I136").Value = 2995.6538
$ws.Range("K136").Value = 8986.9614
$ws.Range("M136").Value = -6436.9614

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 5794
$ws.Range("J125").Value = 5800
$ws.Range("L125").Value = 17400
$ws.Range("N125").Value = -27240
$ws.Range("H128").Value = 144868.5
$ws.Range("I128").Value = 144868.5
$ws.Range("K128").Value = 434605.5
$ws.Range("M128").Value = -429625.5
$ws.Range("H129").Value = 182872.73
$ws.Range("J129").Value = 286895.28
$ws.Range("L129").Value = 860685.8400000001
$ws.Range("N129").Value = -870685.8400000001
$ws.Range("H131").Value = 25069.791
$ws.Range("J131").Value = 30389.885
$ws.Range("L131").Value = 91169.655
$ws.Range("N131").Value = -101249.655
$ws.Range("H137").Value = 77901.07000000001
$ws.Range("I137").Value = 66089.81
$ws.Range("J137").Value = 95081.09
$ws.Range("K137").Value = 198269.43
$ws.Range("L137").Value = 285243.27
$ws.Range("M137").Value = -193169.43
$ws.Range("N137").Value = -295443.27

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 61753.055
$ws.Range("I80").Value = 8073.4443
$ws.Range("J80").Value = 115432.664
$ws.Range("K80").Value = 8073.4443
$ws.Range("L80").Value = 115432.664
$ws.Range("M80").Value = -7075.4443
$ws.Range("N80").Value = -117428.664
$ws.Range("H83").Value = 61753.055
$ws.Range("I83").Value = 8073.4443
$ws.Range("J83").Value = 115432.664
$ws.Range("K83").Value = 40367.2215
$ws.Range("L83").Value = 577163.3200000001
$ws.Range("M83").Value = -35375.2215
$ws.Range("N83").Value = -587147.3200000001
$ws.Range("H97").Value = 2129.4
$ws.Range("I97").Value = 1580.6875
$ws.Range("J97").Value = 4324.25
$ws.Range("K97").Value = 1580.6875
$ws.Range("L97").Value = 4324.25
$ws.Range("M97").Value = -1084.6875
$ws.Range("N97").Value = -5316.25
$ws.Range("H102").Value = 3736.875
$ws.Range("I102").Value = 3985.3572
$ws.Range("K102").Value = 3985.3572
$ws.Range("M102").Value = -2363.3572
$ws.Range("H126").Value = 7591.365
$ws.Range("I126").Value = 5321.304
$ws.Range("K126").Value = 15963.912
$ws.Range("M126").Value = -13493.912
$ws.Range("H132").Value = 5439.095
$ws.Range("I132").Value = 1601.0834
$ws.Range("K132").Value = 4803.2502
$ws.Range("M132").Value = -2273.2502

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1936.5
$ws.Range("I82").Value = 1387.1666
$ws.Range("J82").Value = 2485.8333
$ws.Range("K82").Value = 1387.1666
$ws.Range("L82").Value = 2485.8333
$ws.Range("M82").Value = -1026.1666
$ws.Range("N82").Value = -3207.8333
$ws.Range("H85").Value = 1936.5
$ws.Range("I85").Value = 1387.1666
$ws.Range("J85").Value = 2485.8333
$ws.Range("K85").Value = 1387.1666
$ws.Range("L85").Value = 2485.8333
$ws.Range("M85").Value = -139.1666
$ws.Range("N85").Value = -4981.8333
$ws.Range("H122").Value = 4196.2856
$ws.Range("I122").Value = 3302.05
$ws.Range("K122").Value = 9906.150000000001
$ws.Range("M122").Value = -7456.150000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 41500
$ws.Range("J94").Value = 41500
$ws.Range("L94").Value = 41500
$ws.Range("N94").Value = -43302
$ws.Range("H126").Value = 2007.625
$ws.Range("I126").Value = 1647
$ws.Range("J126").Value = 2127.8333
$ws.Range("K126").Value = 4941
$ws.Range("L126").Value = 6383.499899999999
$ws.Range("M126").Value = -2471
$ws.Range("N126").Value = -11323.4999
$ws.Range("H136").Value = 15306525
$ws.Range("I136").Value = 23811770
$ws.Range("K136").Value = 71435310
$ws.Range("M136").Value = -71432760
